$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ENW006 row (row 7): append the new Jira id and the new verification text
# to the existing "||"-delimited lists stored in B7 (Jira id) and C7 (Description).
$ws.Range("B7").Value = "OPQA-1701||OPQA-3290"
$ws.Range("C7").Value = 'Verify that the "Thanks for your interest in EndNote......" modal displayed when user clicks on the export button when user is signed to facebook account and not having existing steam account||Verify that,after clicking "send to Endnote" Button,user should be able to see the text "Enter you existing account credential (CortellisTM, EndNoteTM Online,InCitesTM, ResearcherID,Thomson InnovationTM, Web of ScienceTM) to link your accounts."'

# The row grew taller to fit the extra sentence.
$ws.Rows(7).RowHeight = 75

# Move the active selection to the edited row (C7), matching the saved view state.
[void]$ws.Range("C7").Select()
